$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data to reflect latest scrape
$ws.Range("D2").Value = "26.203.81"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "1.653.43"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'219.06"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").Value = "'0.5226"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").Value = "'0.2630"
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("D9").Value = "'0.06328"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D10").Value = "'20.46"
$ws.Range("E10").Value = "  -0.68%  "
$ws.Range("D11").Value = "'0.07665"
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("D12").Value = "'4.599"
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("D13").Value = "1.670.98"
$ws.Range("E13").Value = "  +2.02%  "
$ws.Range("D14").Value = "1.876.13"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").Value = "'0.5605"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").Value = "0.0₅8147"
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("D17").Value = "'65.35"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "26.113.53"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("D20").Value = "'4.611"
$ws.Range("E20").Value = "  -2.24%  "
$ws.Range("D21").Value = "'194.95"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("D22").Value = "'10.50"
$ws.Range("E22").Value = "  +2.53%  "
$ws.Range("D23").Value = "'5.948"
$ws.Range("E23").Value = "  -1.18%  "
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").Value = "'145.13"
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("D26").Value = "'0.1189"
$ws.Range("E26").Value = "  -1.37%  "
$ws.Range("D27").Value = "'7.220"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "'15.97"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").Value = "'1.536"
$ws.Range("E29").Value = "  +2.26%  "
$ws.Range("D30").Value = "'0.05488"
$ws.Range("E30").Value = "  -2.05%  "
$ws.Range("D31").Value = "'1.270"
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("D32").Value = "'3.472"
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("D33").Value = "'3.339"
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("D34").Value = "'1.563"
$ws.Range("E34").Value = "  -2.13%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.415"
$ws.Range("E35").Value = "  +0.43%  "
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").Value = "'2.785"
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("D37").Value = "'0.9451"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").Value = "'0.5631"
$ws.Range("E38").Value = "  -1.88%  "
$ws.Range("D39").Value = "'0.01577"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("E40").Value = "  -3.71%  "
$ws.Range("D41").Value = "'1.002"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").Value = "1.029.05"
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("D43").Value = "'0.8193"
$ws.Range("E43").Value = "  -3.22%  "
$ws.Range("D44").Value = "'100.59"
$ws.Range("E44").Value = "  -1.83%  "
$ws.Range("D45").Value = "1.787.66"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'57.48"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₈108"
$ws.Range("E47").Value = "  +1.91%  "
$ws.Range("E48").Value = "  +0.18%  "
$ws.Range("D49").Value = "'0.4334"
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("D50").Value = "'7.886"
$ws.Range("E50").Value = "  -2.00%  "
$ws.Range("D51").Value = "'0.05131"
$ws.Range("E51").Value = "  -3.45%  "

# Ensure column D retains its original "General" cell style/number format
# (plain-text numeric-looking values were entered with a leading apostrophe
# above so Excel doesn't silently convert them to floating point numbers).
$ws.Range("D2:D51").Style = "Normal"
